# Apply updated knapsack-voting ballots (boolean approval matrix, rows 2-22)
# to Sheet1 of the workbook. Only the cells that actually flipped value in
# the source diff are touched; every other cell is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = $false
$ws.Range("C2").Value = $false

$ws.Range("B3").Value = $true
$ws.Range("C3").Value = $true
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("G3").Value = $true

$ws.Range("B4").Value = $false
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false

$ws.Range("B5").Value = $true
$ws.Range("F5").Value = $true
$ws.Range("G5").Value = $true

$ws.Range("E7").Value = $false
$ws.Range("G7").Value = $false

$ws.Range("C8").Value = $true
$ws.Range("D8").Value = $true
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = $true

$ws.Range("C9").Value = $true
$ws.Range("D9").Value = $true
$ws.Range("F9").Value = $true

$ws.Range("B10").Value = $false
$ws.Range("F10").Value = $false

$ws.Range("B12").Value = $true
$ws.Range("C12").Value = $true
$ws.Range("D12").Value = $true
$ws.Range("F12").Value = $true

$ws.Range("B13").Value = $false
$ws.Range("C13").Value = $false
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = $false
$ws.Range("G13").Value = $false

$ws.Range("F15").Value = $false

$ws.Range("B16").Value = $true
$ws.Range("D16").Value = $true

$ws.Range("B17").Value = $false
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = $false
$ws.Range("E17").Value = $false
$ws.Range("F17").Value = $false
$ws.Range("G17").Value = $false

$ws.Range("E18").Value = $false
$ws.Range("G18").Value = $false

$ws.Range("B21").Value = $false
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = $false
$ws.Range("E21").Value = $false
$ws.Range("F21").Value = $false
$ws.Range("G21").Value = $false

$ws.Range("B22").Value = $false
$ws.Range("C22").Value = $false
$ws.Range("D22").Value = $false
$ws.Range("E22").Value = $false
$ws.Range("F22").Value = $false
